$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 386.42856
$ws.Range("I4").Value = 386.42856
$ws.Range("K4").Value = 386.42856
$ws.Range("M4").Value = -272.42856
$ws.Range("H5").Value = 4000601.8
$ws.Range("I5").Value = 2500656
$ws.Range("K5").Value = 2500656
$ws.Range("M5").Value = -2500541
$ws.Range("H19").Value = 887.73334
$ws.Range("I19").Value = 697.6667
$ws.Range("J19").Value = 1014.44446
$ws.Range("K19").Value = 697.6667
$ws.Range("L19").Value = 1014.44446
$ws.Range("M19").Value = -522.6667
$ws.Range("N19").Value = -1364.44446
$ws.Range("H28").Value = 8315.916999999999
$ws.Range("I28").Value = 1665.7778
$ws.Range("J28").Value = 28266.334
$ws.Range("K28").Value = 1665.7778
$ws.Range("L28").Value = 28266.334
$ws.Range("M28").Value = -1180.7778
$ws.Range("N28").Value = -29236.334
$ws.Range("H40").Value = 2755.2222
$ws.Range("I40").Value = 2149.5
$ws.Range("J40").Value = 3966.6667
$ws.Range("K40").Value = 2149.5
$ws.Range("L40").Value = 3966.6667
$ws.Range("M40").Value = -1974.5
$ws.Range("N40").Value = -4316.6667
$ws.Range("H94").Value = 2694.2
$ws.Range("I94").Value = 2694.2
$ws.Range("K94").Value = 2694.2
$ws.Range("M94").Value = -2243.2
$ws.Range("H100").Value = 2457.875
$ws.Range("I100").Value = 2577.5
$ws.Range("J100").Value = 2099
$ws.Range("K100").Value = 2577.5
$ws.Range("L100").Value = 2099
$ws.Range("M100").Value = -2036.5
$ws.Range("N100").Value = -3181
$ws.Range("H107").Value = 4433.5415
$ws.Range("I107").Value = 3970.5
$ws.Range("K107").Value = 3970.5
$ws.Range("M107").Value = -2050.5
$ws.Range("H111").Value = 703
$ws.Range("I111").Value = 393.8889
$ws.Range("J111").Value = 1398.5
$ws.Range("K111").Value = 1181.6667
$ws.Range("L111").Value = 4195.5
$ws.Range("M111").Value = 1885.3333
$ws.Range("N111").Value = -10329.5
$ws.Range("H138").Value = 837689.2
$ws.Range("J138").Value = 951696
$ws.Range("L138").Value = 2855088
$ws.Range("N138").Value = -2865368
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3676.3208
$ws.Range("J32").Value = 11117.286
$ws.Range("L32").Value = 11117.286
$ws.Range("N32").Value = -11691.286
$ws.Range("H61").Value = 4055.484
$ws.Range("I61").Value = 2908.2856
$ws.Range("J61").Value = 6464.6
$ws.Range("K61").Value = 2908.2856
$ws.Range("L61").Value = 6464.6
$ws.Range("M61").Value = -2696.2856
$ws.Range("N61").Value = -6888.6
$ws.Range("H136").Value = 4055.484
$ws.Range("I136").Value = 2908.2856
$ws.Range("J136").Value = 6464.6
$ws.Range("K136").Value = 8724.856800000001
$ws.Range("L136").Value = 19393.8
$ws.Range("M136").Value = -6174.856800000001
$ws.Range("N136").Value = -24493.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 80313.664
$ws.Range("J96").Value = 100471
$ws.Range("L96").Value = 100471
$ws.Range("N96").Value = -105963
$ws.Range("H99").Value = 3868.7856
$ws.Range("I99").Value = 3715.0908
$ws.Range("K99").Value = 3715.0908
$ws.Range("M99").Value = -2217.0908
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1848.7142
$ws.Range("I94").Value = 2107.7334
$ws.Range("K94").Value = 2107.7334
$ws.Range("M94").Value = -1656.7334
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16766889
$ws.Range("I4").Value = 3551954
$ws.Range("K4").Value = 10655862
$ws.Range("M4").Value = -10655750
$ws.Range("H13").Value = 999
$ws.Range("I13").Value = 999
$ws.Range("K13").Value = 2997
$ws.Range("M13").Value = -2829
$ws.Range("H98").Value = 448.33334
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 448.33334
$ws.Range("K98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("M98").Value = 1345.00002
$ws.Range("N98").Value = -4341.000019999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 12000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 12000
$ws.Range("K41").Value = 0
$ws.Range("L41").ClearContents()
$ws.Range("M41").Value = 12000
$ws.Range("N41").Value = -12710
$ws.Range("H70").Value = 11360
$ws.Range("J70").Value = 11360
$ws.Range("L70").Value = 11360
$ws.Range("N70").Value = -11900
$ws.Range("H73").Value = 11360
$ws.Range("J73").Value = 11360
$ws.Range("L73").Value = 11360
$ws.Range("N73").Value = -13232
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 0
$ws.Range("H107").Value = 650.2857
$ws.Range("I107").Value = 187.875
$ws.Range("J107").Value = 1266.8334
$ws.Range("K107").Value = 187.875
$ws.Range("L107").Value = 1266.8334
$ws.Range("M107").Value = 1732.125
$ws.Range("N107").Value = -5106.8334
$ws.Range("H126").Value = 4770.9443
$ws.Range("I126").Value = 3784.875
$ws.Range("J126").Value = 5559.8
$ws.Range("K126").Value = 11354.625
$ws.Range("L126").Value = 16679.4
$ws.Range("M126").Value = -8884.625
$ws.Range("N126").Value = -21619.4
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9615.25
$ws.Range("I7").Value = 12321
$ws.Range("J7").Value = 1498
$ws.Range("K7").Value = 12321
$ws.Range("L7").Value = 1498
$ws.Range("M7").Value = -12209
$ws.Range("N7").Value = -1722
$ws.Range("H40").Value = 2345.75
$ws.Range("I40").Value = 2109.4285
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 2109.4285
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -1973.4285
$ws.Range("N40").Value = -4272
$ws.Range("H43").Value = 5967333.5
$ws.Range("I43").Value = 3139500
$ws.Range("K43").Value = 3139500
$ws.Range("M43").Value = -3139307
$ws.Range("H100").Value = 6824
$ws.Range("I100").Value = 6792
$ws.Range("K100").Value = 6792
$ws.Range("M100").Value = -6251
$ws.Range("H122").Value = 4833
$ws.Range("I122").Value = 3999.5
$ws.Range("K122").Value = 11998.5
$ws.Range("M122").Value = -9548.5
$ws.Range("H126").Value = 9615.25
$ws.Range("I126").Value = 12321
$ws.Range("J126").Value = 1498
$ws.Range("K126").Value = 36963
$ws.Range("L126").Value = 4494
$ws.Range("M126").Value = -34493
$ws.Range("N126").Value = -9434
$ws.Range("H132").Value = 6732.9375
$ws.Range("J132").Value = 8959.625
$ws.Range("L132").Value = 26878.875
$ws.Range("N132").Value = -31938.875
$ws.Range("H136").Value = 5937.476
$ws.Range("I136").Value = 3461.7693
$ws.Range("K136").Value = 10385.3079
$ws.Range("M136").Value = -7835.3079
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 27499.5
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H96").Value = 2023319.6
$ws.Range("I96").Value = 4043873.8
$ws.Range("K96").Value = 4043873.8
$ws.Range("M96").Value = -4042500.8
$ws.Range("H132").Value = 5758
$ws.Range("I132").Value = 5585.857
$ws.Range("K132").Value = 16757.571
$ws.Range("M132").Value = -14227.571
$ws.Range("H136").Value = 3151.2632
$ws.Range("J136").Value = 1850
$ws.Range("L136").Value = 5550
$ws.Range("N136").Value = -10650
